$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing "no resuelto" / "No resuelto" annotations from the
# comments in rows 4, 5 and 6 (issue now fixed, note removed).
$ws.Range("G4").Value = "Recuerda que si el cliente despues de tres intentos no contesta agregar esa información"
$ws.Range("G5").Value = "En plan de riesgos falta efecto de la causa no conexión por admiadmin impide conectar a maquina cliente por ejemplo."
$ws.Range("G6").Value = "Al no requerir implementacion no debe decir los mensajes de implementacion"

# Row 7 (id 4, "En equipo de trabajo agregar nombre de empresa") is now
# resolved: status flips to "Cerrada" and the comment drops its
# "no resuelto" suffix.
$ws.Range("F7").Value = "Cerrada"
$ws.Range("G7").Value = "En la parte equipo de empresa poner SOS Software y en cliente el nombre de cliente"

# Row 8 (id 5, "Cambiar nombre pland e proyecto a plan_de_proyecto") is
# merged away entirely, leaving only a stray "cd " text in A8.
$ws.Range("B8:G8").ClearContents()
$ws.Range("A8").Value = "cd "

# Row heights shrink to match the shorter wrapped text.
$ws.Rows.Item(5).RowHeight = 55.2
$ws.Rows.Item(7).RowHeight = 41.75
$ws.Rows.Item(8).RowHeight = 13.8

$ws.Range("C8").Select()
